# Insert a new weekly price-report row at row 15 (shifts existing rows 15-139
# down to 16-140) and populate it with the new record's data, per the
# commit "Fruta / hortaliza, semanal".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 15; this shifts rows 15..139 down
# to 16..140 (and carries the D-column's date style down with them).
$ws.Rows.Item(15).Insert()

# Columns that are constant for every data row in this sheet.
$ws.Range("A15").Value = 8
$ws.Range("B15").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C15").Value = 'Coquimbo'
$ws.Range("E15").Value = 4
$ws.Range("F15").Value = 100112052
$ws.Range("G15").Value = 'Albahaca'
$ws.Range("H15").Value = 'Sin especificar'
$ws.Range("I15").Value = 'Primera'
$ws.Range("R15").Value = 'Hortaliza'

# New record's own data.
$ws.Range("D15").Value2 = 44901
$ws.Range("J15").Value = 1000
$ws.Range("K15").Value = 4000
$ws.Range("L15").Value = 4500
$ws.Range("M15").Value = 4250
$ws.Range("N15").Value = '$/paquete'
$ws.Range("O15").Value = 'Región de Arica y Parinacota'
$ws.Range("P15").Value = 4250
$ws.Range("Q15").Value = 1
